$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.438.59"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.562.96"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  +1.16%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.499"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  -0.11%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "21.86"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.786.04"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "1.565.03"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("E15").Value = "  -1.88%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "63.33"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "27.423.65"
$ws.Range("E17").Value = "  -0.28%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "213.32"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  -0.11%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "9.55"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  +1.66%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "152.92"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -0.08%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "6.73"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D33").Value = "1.359.92"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  +1.68%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.972"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("E39").Value = "  -0.48%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.821"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("E41").Value = "  -0.13%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.974"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  +1.46%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "64.10"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +1.40%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "5.28"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "1.699.86"
$ws.Range("E47").Value = "  -0.57%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "85.38"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("D49").Value = "0.0₇0989"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("E51").Value = "  -0.43%  "
